$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header label change
$ws.Range("O1").Value = "F1 train"

# Row 2 (SVM)
$ws.Range("O2").Value = 1

# Row 3 (RF)
$ws.Range("O3").Value = 0.9117647058823529

# Row 4 (XGB)
$ws.Range("O4").Value = 1

# Row 5 (KNN)
$ws.Range("O5").Value = 0.6984126984126984

# Row 6 (MLP) - params + confusion matrix + metrics
$ws.Range("C6").Value = "{'activation': 'relu', 'alpha': 0.0001, 'hidden_layer_sizes': (64, 32), 'learning_rate': 'constant'}"
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 8
$ws.Range("I6").Value = 0.7
$ws.Range("J6").Value = 0.6666666666666666
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.7272727272727273
$ws.Range("O6").Value = 0.7714285714285715

# Row 8 (RF, 10%)
$ws.Range("O8").Value = 1

# Row 9 (XGB, 10%)
$ws.Range("O9").Value = 1

# Row 10 (KNN, 10%)
$ws.Range("O10").Value = 1

# Row 11 (MLP, 10%) - params + confusion matrix + metrics
$ws.Range("C11").Value = "{'activation': 'tanh', 'alpha': 0.0001, 'hidden_layer_sizes': (64,), 'learning_rate': 'constant'}"
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 9
$ws.Range("H11").Value = 5
$ws.Range("I11").Value = 0.65
$ws.Range("J11").Value = 0.5333333333333333
$ws.Range("K11").Value = 0.4444444444444444
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.8181818181818182
$ws.Range("N11").Value = 0.4444444444444444
$ws.Range("O11").Value = 0.704225352112676

# Row 13 (RF, Free)
$ws.Range("O13").Value = 0.9714285714285714

# Row 14 (XGB, Free)
$ws.Range("O14").Value = 0.9859154929577465

# Row 15 (KNN, Free)
$ws.Range("O15").Value = 1

# Row 16 (MLP, Free)
$ws.Range("O16").Value = 0.6363636363636364
